# Generate Report for Handoff
# The "eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3" file moved from "Handed back: in
# sync with en-US" to "Ready for handoff" with an updated handoff timestamp
# and a new error detail message noting the handback file is stale.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet — row 3 is the eaeb81fe-... entry.
#   E3 = zh-cn status, F3 = de-de status, G3 = Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-04 16:53:30"

# ---------------------------------------------------------------------------
# zh-cn sheet — row 3 is the eaeb81fe-... entry.
#   C3 = Status, H3 = Latest Handoff Datetime, P3 = Error Detail
# ---------------------------------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-04 16:53:25"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2499c35ce336dbca9590f680a136302312995a16/e2e/eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0eb824914866809b1fcf6251deac64c626db1e87/e2e/eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md."
# Column P (16th) widened to fit the new long error message.
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------------
# de-de sheet — row 3 is the eaeb81fe-... entry.
#   C3 = Status, H3 = Latest Handback DateTime, P3 = Error Detail
# ---------------------------------------------------------------------------
$dede = $wb.Sheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-04 16:53:30"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2499c35ce336dbca9590f680a136302312995a16/e2e/eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0eb824914866809b1fcf6251deac64c626db1e87/e2e/eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md."
# Column P (16th) widened to fit the new long error message.
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
